# Apply "config added removed +1 for threshold in config" edit:
# Change the H/I column formulas from 1/(SUM(Fx,1)) and 1/(SUM(Gx,1))
# (which effectively added +1 to the denominator) to IFERROR(1/Fx,0)
# and IFERROR(1/Gx,0) respectively, removing that +1 threshold offset
# and guarding against division by zero with IFERROR.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds standalone (non-shared) formulas.
$ws.Range("H2").Formula = "=IFERROR(1/F2,0)"
$ws.Range("I2").Formula = "=IFERROR(1/G2,0)"

# Rows 3:25 form the shared-formula blocks anchored at H3/I3.
$ws.Range("H3:H25").Formula = "=IFERROR(1/F3,0)"
$ws.Range("I3:I25").Formula = "=IFERROR(1/G3,0)"

# Update the saved cursor/selection position to I25.
$ws.Range("I25").Select()
